# "stories done, game mechanics"
#
# The tracked progress sheet gets two data updates:
#   - B4 ("Theory" Written pages):  4  -> 8
#   - I12 (pages written on the day tracked in row 12): 0 -> 4
#
# Every other cell that differs in the target workbook (B11, D11, B12,
# I31, and the TODAY()-driven F2/F3/F4 trio) is a formula that recalculates
# automatically from these two inputs, so we only need to touch the two
# literal values below and let Excel's recalculation do the rest.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 8
$ws.Range("I12").Value = 4
